$d = $word.ActiveDocument

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# -----------------------------------------------------------------
# Edit 1: insert a new paragraph before "Operaciones o funcionalidades:"
# -----------------------------------------------------------------
$target1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Operaciones o funcionalidades:*") {
        $target1 = $p
        break
    }
}

if ($target1 -ne $null) {
    $r = $target1.Range.Duplicate
    $r.Collapse(1)
    $r.InsertParagraphBefore()
    $newR = $r.Duplicate
    $newR.Collapse(1)

    $body1 = '<w:p><w:pPr><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">Una duda que me surgi' + [char]0x00F3 + ' en este </w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>punto,</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> es si la estructura de paquetes que defin' + [char]0x00ED + ' para la capa de aplicaci' + [char]0x00F3 + 'n era correcta. Una vez creado el primer puerto de entrada y la primera clase de datos encapsulados para su uso, me di cuenta de que necesitar' + [char]0x00ED + 'a servicios para registrar el usuario (quiz' + [char]0x00E1 + ' tambi' + [char]0x00E9 + 'n repositorios para guardarlos), la soluci' + [char]0x00F3 + 'n fue preguntarle a la IA que deber' + [char]0x00ED + 'a hacer.</w:t></w:r>' + `
        '</w:p>'

    $newR.InsertXML($pkgHeader + $body1 + $pkgFooter)
}

# -----------------------------------------------------------------
# Edit 2: insert four new list paragraphs after
# "Al ser una arquitectura hexagonal..."
# -----------------------------------------------------------------
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Al ser una arquitectura hexagonal*") {
        $target2 = $p
    }
}

if ($target2 -ne $null) {
    $nextPara2 = $target2.Next()
    $r2 = $nextPara2.Range.Duplicate
    $r2.Collapse(1)
    $r2.InsertParagraphBefore()
    $newR2 = $r2.Duplicate
    $newR2.Collapse(1)

    $para1 = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">Por otro lado, empec' + [char]0x00E9 + ' haciendo los puertos con </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>DTOs</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">, pero seg' + [char]0x00FA + 'n las buenas pr' + [char]0x00E1 + 'cticas de arquitectura hexagonal y DDD lo correcto ser' + [char]0x00ED + 'a usar </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>Commands</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>:</w:t></w:r>' + `
        '</w:p>'

    $para2 = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>Encapsula los datos de entrada.</w:t></w:r>' + `
        '</w:p>'

    $para3 = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">Permite hacer validaciones previas al uso. </w:t></w:r>' + `
        '</w:p>'

    $para4 = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">Logramos desacoplarlo de la interfaz. </w:t></w:r>' + `
        '</w:p>'

    $para5 = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">Cliff utilizo </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>Commands</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> varias veces en las clases.</w:t></w:r>' + `
        '</w:p>'

    $body2 = $para1 + $para2 + $para3 + $para4 + $para5

    $newR2.InsertXML($pkgHeader + $body2 + $pkgFooter)
}
